# Minha primeira versão com CustomTkinter
# Extends the "Registros" (attendance) sheet with extra date columns
# (Q:X) for November 2025 plus the matching attendance marks for the
# existing students in rows 3-5 (row 2 already had values in H:P).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registros")

# ---------------------------------------------------------------------
# 1) New header dates in row 1 (Q1:X1), bold / centered / bordered like
#    the rest of the header row (same visual style as H1:P1).
#    NumberFormat is forced to Text first so the dd/mm/yyyy-looking
#    strings are stored as literal text instead of being auto-parsed
#    into date serial numbers; the formats are then re-applied from an
#    existing header cell (P1) so the final look matches the header.
# ---------------------------------------------------------------------
$headerDates = @("04/11/2025","06/11/2025","11/11/2025","13/11/2025","18/11/2025","20/11/2025","25/11/2025","27/11/2025")
$headerCols = @(17,18,19,20,21,22,23,24)   # Q..X

$ws.Range("Q1:X1").NumberFormat = "@"
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Cells.Item(1, $headerCols[$i]).Value = $headerDates[$i]
}

$ws.Range("P1").Copy()
$ws.Range("Q1:X1").PasteSpecial(-4122)   # xlPasteFormats - reapply header look (bold, border, centered)

# ---------------------------------------------------------------------
# 2) Fill in the previously-blank attendance marks for "zé do pé" (row 3),
#    "tonho do sonho" (row 4) and "ana banana" (row 5) in columns H, I, J,
#    K, M, and populate the new November columns Q:X for every student row.
# ---------------------------------------------------------------------
$cols = @(8, 9, 10, 11, 13, 17, 18, 19, 20, 21, 22, 23, 24)   # H, I, J, K, M, Q, R, S, T, U, V, W, X

$row3 = @("j", "f", "c", "c", "c", "c", "f", "f", "j", "j", "j", "c", "c")
$row4 = @("j", "f", "c", "f", "f", "c", "c", "c", "c", "c", "j", "c", "c")
$row5 = @("j", "c", "c", "j", "j", "c", "c", "c", "c", "c", "j", "c", "c")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(3, $cols[$i]).Value = $row3[$i]
    $ws.Cells.Item(4, $cols[$i]).Value = $row4[$i]
    $ws.Cells.Item(5, $cols[$i]).Value = $row5[$i]
}
